# Update the "取得日時" (acquisition datetime) column (A) for all data rows
# on the "ランサーズ" sheet from "2025-09-10 12:40:01" to "2025-09-10 12:53:42".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-09-10 12:40:01"
$newValue = "2025-09-10 12:53:42"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 2
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
